$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update user/password values in row 2 and row 3
$ws.Range("A2").Value = "james8928748234"
$ws.Range("B3").Value = "Tommy82379834893"
$ws.Range("A3").Value = "thomas798597241"

# Adjust column widths (target OOXML widths: 16.88671875 and 20.6640625;
# the nearest values achievable through the ColumnWidth property's
# internal rounding are used here)
$ws.Columns.Item(1).ColumnWidth = 16
$ws.Columns.Item(2).ColumnWidth = 19.833333333333336

# Update selection to A3
$ws.Range("A3").Select()
